$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 ("area" / " QUALITY") and shift rows 6-8 up.
$ws.Rows("5:5").Delete()
